# Add files via upload
# The header row previously used "<Name>-Norm" labels; rename them back to
# the plain column names, and move the active-cell selection from H10 to H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = "Estacion"
$ws.Range("B1").Value = "Prof"
$ws.Range("C1").Value = "PesoEsp"
$ws.Range("D1").Value = "RCS"
$ws.Range("E1").Value = "RQD"
$ws.Range("F1").Value = "RMR"
$ws.Range("G1").Value = "GSI"
$ws.Range("H1").Value = "TipoFort"

$ws.Activate()
$ws.Range("H2").Select()
